$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$text)
    # Force the cell to store a literal string even when the text looks like
    # a number (document IDs, phone numbers, passwords) - matches how the
    # workbook already stores this kind of data as shared strings instead of
    # numeric values - then drop the temporary "@" text format so the cell
    # keeps the workbook's default style.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- "docentes" sheet (sheet1): add new row 4 with the "andres lopez" record ---
$wsDocentes = $wb.Worksheets.Item("docentes")

$wsDocentes.Range("A4").Value = "andres"
$wsDocentes.Range("B4").Value = "lopez"
Set-TextValue $wsDocentes.Range("C4") "1233445678"
$wsDocentes.Range("D4").Value = "andres01"
$wsDocentes.Range("E4").Value = "andres@gmail.com"
Set-TextValue $wsDocentes.Range("F4") "3001234567"
Set-TextValue $wsDocentes.Range("G4") "12345"
$wsDocentes.Range("H4").Value = "Docente"
$wsDocentes.Range("I4").Value = "Antioquia"
$wsDocentes.Range("J4").Value = "Medellín"
$wsDocentes.Range("K4").Value = "P.C.J.I.C"

# --- "estudiantes" sheet (sheet2): add new row 5 with the same person, registered as student ---
$wsEstudiantes = $wb.Worksheets.Item("estudiantes")

$wsEstudiantes.Range("A5").Value = "andres"
$wsEstudiantes.Range("B5").Value = "lopez"
Set-TextValue $wsEstudiantes.Range("C5") "1234567789"
$wsEstudiantes.Range("D5").Value = "andres01"
$wsEstudiantes.Range("E5").Value = "andres@gmail.com"
Set-TextValue $wsEstudiantes.Range("F5") "3015416963"
Set-TextValue $wsEstudiantes.Range("G5") "12345"
$wsEstudiantes.Range("H5").Value = "Estudiante"
$wsEstudiantes.Range("I5").Value = "Antioquia"
$wsEstudiantes.Range("J5").Value = "Medellín"
$wsEstudiantes.Range("K5").Value = "P.C.J.I.C"

# --- Switch the active tab from "docentes" to "estudiantes" (matches activeTab/tabSelected change) ---
$wsEstudiantes.Activate()
